# Horarios actualizados Linea 141 - 249
# New scrape timestamp / updated data rows.

$wb = $excel.ActiveWorkbook

$newTime = "03:20:50"

# ---------------------------------------------------------------------
# Sheet "LP1912": had 3 data rows (6,7,8) -> now only 1 data row (6),
# which now carries the old row-8 values (215A_EL PATO) with updated
# Hora_Llegada / Minutos.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 1"

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "04:45"
$ws1.Range("C6").Value = "215A_EL PATO"
$ws1.Range("D6").Value = 85
$ws1.Range("E6").Value = "LP1912"

# Drop the now-obsolete rows 7 and 8 (shifts dimension to A1:E6).
$ws1.Rows("7:8").Delete()

# ---------------------------------------------------------------------
# Sheet "LP1912-215": had 2 data rows (6,7) -> now only 1 data row (6),
# which now carries the old row-7 values (215A_EL PATO) with updated
# Hora_Llegada / Minutos.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "04:45"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 85
$ws2.Range("E6").Value = "LP1912"

$ws2.Rows("7:7").Delete()

# ---------------------------------------------------------------------
# Sheet "6203-6173": only the "Última actualización" timestamp moves;
# "Total filas" stays at 0 (no data rows either before or after).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
